$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '69.279.44'
Set-TextValue 'E2' '  +2.28%  '

Set-TextValue 'D3' '3.397.75'
Set-TextValue 'E3' '  +2.16%  '

Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.10%  '

Set-TextValue 'D5' '587.77'
Set-TextValue 'E5' '  +0.85%  '

Set-TextValue 'D6' '180.93'
Set-TextValue 'E6' '  +3.76%  '

Set-TextValue 'E7' '  +1.77%  '

Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  +0.01%  '

Set-TextValue 'D9' '0.202'
Set-TextValue 'E9' '  +11.35%  '

Set-TextValue 'E10' '  +2.68%  '

Set-TextValue 'D11' '48.45'
Set-TextValue 'E11' '  +3.26%  '

Set-TextValue 'D12' '0.0000287'
Set-TextValue 'E12' '  +5.64%  '

Set-TextValue 'D13' '683.93'
Set-TextValue 'E13' '  -1.83%  '

Set-TextValue 'E14' '  +4.16%  '

Set-TextValue 'D15' '3.946.38'
Set-TextValue 'E15' '  +2.07%  '

Set-TextValue 'D16' '69.344.61'
Set-TextValue 'E16' '  +2.29%  '

Set-TextValue 'B17' 'TRON'
Set-TextValue 'C17' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D17' '0.121'
Set-TextValue 'E17' '  +1.62%  '

Set-TextValue 'B18' 'WrappedEther'
Set-TextValue 'C18' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '3.392.05'
Set-TextValue 'E18' '  +2.76%  '

Set-TextValue 'D19' '17.73'
Set-TextValue 'E19' '  +1.67%  '

Set-TextValue 'E20' '  +2.01%  '

Set-TextValue 'D21' '0.912'
Set-TextValue 'E21' '  +2.65%  '

Set-TextValue 'D22' '17.32'
Set-TextValue 'E22' '  +2.52%  '

Set-TextValue 'D23' '5.36'
Set-TextValue 'E23' '  -0.98%  '

Set-TextValue 'D24' '103.31'

Set-TextValue 'D25' '3.94'
Set-TextValue 'E25' '  +1.15%  '

Set-TextValue 'D26' '2.73'
Set-TextValue 'E26' '  +2.03%  '

Set-TextValue 'D27' '9.68'
Set-TextValue 'E27' '  +2.96%  '

Set-TextValue 'D28' '33.90'
Set-TextValue 'E28' '  +3.29%  '

Set-TextValue 'E29' '  +3.49%  '

Set-TextValue 'D30' '6.96'
Set-TextValue 'E30' '  -0.05%  '

Set-TextValue 'E31' '  +1.64%  '

Set-TextValue 'D32' '557.00'
Set-TextValue 'E32' '  -2.85%  '

Set-TextValue 'D33' '3.62'
Set-TextValue 'E33' '  +11.03%  '

Set-TextValue 'E34' '  +1.62%  '

Set-TextValue 'D35' '58.63'
Set-TextValue 'E35' '  +3.86%  '

Set-TextValue 'E36' '  -0.05%  '

Set-TextValue 'D37' '3.657.25'
Set-TextValue 'E37' '  -1.79%  '

Set-TextValue 'E38' '  +6.23%  '

Set-TextValue 'D39' '36.02'
Set-TextValue 'E39' '  +1.40%  '

Set-TextValue 'D40' '0.0₃0721'
Set-TextValue 'E40' '  +7.72%  '

Set-TextValue 'D41' '3.25'
Set-TextValue 'E41' '  +3.87%  '

Set-TextValue 'D42' '2.67'
Set-TextValue 'E42' '  +2.75%  '

Set-TextValue 'B43' 'VeChain'
Set-TextValue 'C43' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D43' '0.0427'
Set-TextValue 'E43' '  +5.57%  '

Set-TextValue 'B44' 'TheGraph'
Set-TextValue 'C44' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D44' '0.339'
Set-TextValue 'E44' '  +1.62%  '

Set-TextValue 'B45' 'ThetaToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D45' '2.67'
Set-TextValue 'E45' '  +1.82%  '

Set-TextValue 'B46' 'Stellar'
Set-TextValue 'C46' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D46' '0.130'
Set-TextValue 'E46' '  +1.09%  '

Set-TextValue 'B47' 'Mantle'
Set-TextValue 'C47' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D47' '1.39'
Set-TextValue 'E47' '  +5.21%  '

Set-TextValue 'B48' 'FirstDigitalUSD'
Set-TextValue 'C48' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D48' '1.00'
Set-TextValue 'E48' '  -0.06%  '

Set-TextValue 'B49' 'Monero'
Set-TextValue 'C49' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D49' '131.42'
Set-TextValue 'E49' '  +0.34%  '

Set-TextValue 'B50' 'CoreDAO'
Set-TextValue 'C50' 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-TextValue 'D50' '2.73'
Set-TextValue 'E50' '  +3.36%  '

Set-TextValue 'B51' 'THORChain'
Set-TextValue 'C51' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D51' '7.50'
Set-TextValue 'E51' '  +2.35%  '
